$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Disable ("N") the Runmode flag for all test cases except the first two
# (rows 4 through 26), leaving rows 2 and 3 ("Y") enabled.
$ws.Range("D4:D26").Value = "N"

$ws.Activate()
$ws.Range("D6").Select()
